$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 137
$ws.Range("B1").Value = 183.3999999999996
$ws.Range("C1").Value = 137

$ws.Range("A2").Value = 137
$ws.Range("B2").Value = 137
$ws.Range("C2").Value = 137
